# Ansprache wird jetzt unterstuetzt: add ", Hand" after "Psalm23" in the
# document's default (primary) header.
$d = $word.ActiveDocument

$hdr = $d.Sections(1).Headers(1)  # wdHeaderFooterPrimary

$found = $hdr.Range.Find.Execute("Psalm23", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "Psalm23, Hand", 2)

if (-not $found) {
    # Fallback: locate whichever header/footer actually holds "Psalm23"
    # and append the addition directly after the existing text run.
    foreach ($sec in $d.Sections) {
        foreach ($h in $sec.Headers) {
            if ($h.Exists -and ($h.Range.Text -like "*Psalm23*")) {
                $p = $h.Range.Paragraphs(1)
                $r = $p.Range
                $r.SetRange($r.End - 1, $r.End - 1)
                $r.InsertAfter(", Hand")
            }
        }
    }
}
